$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header J1 = PEDC
$ws.Range("J1").Value = "PEDC"

# Clear WESM (H) data column values for rows 2-25 (column header H1 stays)
$ws.Range("H2:H25").ClearContents()

# Update data rows 2-25 with new values per column
# Row 2
$ws.Range("A2").Value = 49055.70767884596
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 12500
$ws.Range("D2").Value = 10000
$ws.Range("E2").Value = 5000
$ws.Range("F2").Value = 5000
$ws.Range("G2").Value = 20000
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5000

# Row 3
$ws.Range("A3").Value = 45913.44015241734
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 12500
$ws.Range("D3").Value = 10000
$ws.Range("E3").Value = 5000
$ws.Range("F3").Value = 5000
$ws.Range("G3").Value = 10000
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 5000

# Row 4
$ws.Range("A4").Value = 43785.45451533525
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = 12500
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 5000
$ws.Range("F4").Value = 5000
$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 5000

# Row 5
$ws.Range("A5").Value = 41929.94255221684
$ws.Range("B5").Value = 5000
$ws.Range("C5").Value = 12500
$ws.Range("D5").Value = 10000
$ws.Range("E5").Value = 5000
$ws.Range("F5").Value = 5000
$ws.Range("G5").Value = 0
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 5000

# Row 6
$ws.Range("A6").Value = 42771.98937702928
$ws.Range("B6").Value = 5000
$ws.Range("C6").Value = 12500
$ws.Range("D6").Value = 10000
$ws.Range("E6").Value = 5000
$ws.Range("F6").Value = 5000
$ws.Range("G6").Value = 0
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 5000

# Row 7
$ws.Range("A7").Value = 43594.49335046268
$ws.Range("B7").Value = 5000
$ws.Range("C7").Value = 12500
$ws.Range("D7").Value = 10000
$ws.Range("E7").Value = 5000
$ws.Range("F7").Value = 5000
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 5000

# Row 8
$ws.Range("A8").Value = 42603.55084781715
$ws.Range("B8").Value = 5000
$ws.Range("C8").Value = 12500
$ws.Range("D8").Value = 10000
$ws.Range("E8").Value = 5000
$ws.Range("F8").Value = 5000
$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 5000

# Row 9
$ws.Range("A9").Value = 35409.32129739954
$ws.Range("B9").Value = 5000
$ws.Range("C9").Value = 12500
$ws.Range("D9").Value = 10000
$ws.Range("E9").Value = 5000
$ws.Range("F9").Value = 5000
$ws.Range("G9").Value = 0
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 5000

# Row 10
$ws.Range("A10").Value = 57208.2105
$ws.Range("B10").Value = 5000
$ws.Range("C10").Value = 25000
$ws.Range("D10").Value = 20000
$ws.Range("E10").Value = 10000
$ws.Range("F10").Value = 10000
$ws.Range("G10").Value = 20000
$ws.Range("I10").Value = 9
$ws.Range("J10").Value = 5000

# Row 11
$ws.Range("A11").Value = 59820.862
$ws.Range("B11").Value = 5000
$ws.Range("C11").Value = 25000
$ws.Range("D11").Value = 20000
$ws.Range("E11").Value = 10000
$ws.Range("F11").Value = 10000
$ws.Range("G11").Value = 20000
$ws.Range("I11").Value = 10
$ws.Range("J11").Value = 5000

# Row 12
$ws.Range("A12").Value = 62649.0755
$ws.Range("B12").Value = 10000
$ws.Range("C12").Value = 25000
$ws.Range("D12").Value = 20000
$ws.Range("E12").Value = 10000
$ws.Range("F12").Value = 10000
$ws.Range("G12").Value = 20000
$ws.Range("I12").Value = 11
$ws.Range("J12").Value = 10000

# Row 13
$ws.Range("A13").Value = 63585.374
$ws.Range("B13").Value = 10000
$ws.Range("C13").Value = 25000
$ws.Range("D13").Value = 20000
$ws.Range("E13").Value = 10000
$ws.Range("F13").Value = 10000
$ws.Range("G13").Value = 20000
$ws.Range("I13").Value = 12
$ws.Range("J13").Value = 10000

# Row 14
$ws.Range("A14").Value = 64148.74625
$ws.Range("B14").Value = 10000
$ws.Range("C14").Value = 25000
$ws.Range("D14").Value = 20000
$ws.Range("E14").Value = 10000
$ws.Range("F14").Value = 10000
$ws.Range("G14").Value = 20000
$ws.Range("I14").Value = 13
$ws.Range("J14").Value = 10000

# Row 15
$ws.Range("A15").Value = 67428.30775000001
$ws.Range("B15").Value = 10000
$ws.Range("C15").Value = 25000
$ws.Range("D15").Value = 20000
$ws.Range("E15").Value = 10000
$ws.Range("F15").Value = 10000
$ws.Range("G15").Value = 20000
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 10000

# Row 16
$ws.Range("A16").Value = 36096.26075
$ws.Range("B16").Value = 10000
$ws.Range("C16").Value = 25000
$ws.Range("D16").Value = 20000
$ws.Range("E16").Value = 10000
$ws.Range("F16").Value = 10000
$ws.Range("G16").Value = 20000
$ws.Range("I16").Value = 15
$ws.Range("J16").Value = 10000

# Row 17
$ws.Range("A17").Value = 0
$ws.Range("B17").Value = 10000
$ws.Range("C17").Value = 25000
$ws.Range("D17").Value = 20000
$ws.Range("E17").Value = 10000
$ws.Range("F17").Value = 10000
$ws.Range("G17").Value = 20000
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 10000

# Row 18
$ws.Range("A18").Value = 1768.269410000001
$ws.Range("B18").Value = 10000
$ws.Range("C18").Value = 25000
$ws.Range("D18").Value = 20000
$ws.Range("E18").Value = 10000
$ws.Range("F18").Value = 10000
$ws.Range("G18").Value = 20000
$ws.Range("I18").Value = 17
$ws.Range("J18").Value = 10000

# Row 19
$ws.Range("A19").Value = 0
$ws.Range("B19").Value = 10000
$ws.Range("C19").Value = 25000
$ws.Range("D19").Value = 20000
$ws.Range("E19").Value = 10000
$ws.Range("F19").Value = 10000
$ws.Range("G19").Value = 20000
$ws.Range("I19").Value = 18
$ws.Range("J19").Value = 10000

# Row 20
$ws.Range("A20").Value = 0
$ws.Range("B20").Value = 10000
$ws.Range("C20").Value = 25000
$ws.Range("D20").Value = 20000
$ws.Range("E20").Value = 10000
$ws.Range("F20").Value = 10000
$ws.Range("G20").Value = 20000
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 10000

# Row 21
$ws.Range("A21").Value = 0
$ws.Range("B21").Value = 10000
$ws.Range("C21").Value = 25000
$ws.Range("D21").Value = 20000
$ws.Range("E21").Value = 10000
$ws.Range("F21").Value = 10000
$ws.Range("G21").Value = 20000
$ws.Range("I21").Value = 20
$ws.Range("J21").Value = 10000

# Row 22
$ws.Range("A22").Value = 0
$ws.Range("B22").Value = 10000
$ws.Range("C22").Value = 25000
$ws.Range("D22").Value = 20000
$ws.Range("E22").Value = 10000
$ws.Range("F22").Value = 10000
$ws.Range("G22").Value = 20000
$ws.Range("I22").Value = 21
$ws.Range("J22").Value = 10000

# Row 23
$ws.Range("A23").Value = 0
$ws.Range("B23").Value = 10000
$ws.Range("C23").Value = 25000
$ws.Range("D23").Value = 20000
$ws.Range("E23").Value = 10000
$ws.Range("F23").Value = 10000
$ws.Range("G23").Value = 20000
$ws.Range("I23").Value = 22
$ws.Range("J23").Value = 10000

# Row 24
$ws.Range("A24").Value = 0
$ws.Range("B24").Value = 10000
$ws.Range("C24").Value = 25000
$ws.Range("D24").Value = 20000
$ws.Range("E24").Value = 10000
$ws.Range("F24").Value = 10000
$ws.Range("G24").Value = 20000
$ws.Range("I24").Value = 23
$ws.Range("J24").Value = 10000

# Row 25
$ws.Range("A25").Value = 0
$ws.Range("B25").Value = 5000
$ws.Range("C25").Value = 12500
$ws.Range("D25").Value = 10000
$ws.Range("E25").Value = 5000
$ws.Range("F25").Value = 5000
$ws.Range("G25").Value = 20000
$ws.Range("I25").Value = 24
$ws.Range("J25").Value = 5000
